$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 202 (SEPIOFF / 1-RAP record), shifting all rows below it
# up by one. This matches the target diff where row 202's old contents
# disappear and every subsequent row's data moves up by one row, with the
# final row (213) vanishing as a result.
$ws.Rows.Item(202).Delete()
